$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title (row 1 keeps its existing merged / coloured style - only the text changes) ---
$ws.Range("A1").Value = "INVENTORY MANAGEMENT SYSTEM"

# --- New header row describing the inventory table columns ---
$ws.Range("A2").Value = "SR.NO"
$ws.Range("B2").Value = "Product Name"
$ws.Range("C2").Value = "Product Id"
$ws.Range("D2").Value = "Product Price"
$ws.Range("E2").Value = "Product Quantity"

# --- Size the columns that hold the inventory table ---
$ws.Range("A:A").ColumnWidth = 10
$ws.Range("B:B").ColumnWidth = 50
$ws.Range("C:C").ColumnWidth = 15
$ws.Range("D:D").ColumnWidth = 15
$ws.Range("E:E").ColumnWidth = 15

# --- Register the bold 14pt font used for the sheet's headings/design ---
$fontCell = $ws.Cells.Item(2, 100)
$fontCell.Font.Bold = $true
$fontCell.Font.Size = 14
$fontCell.Clear()
